$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text for card 4 (row 5): "I couldn't care less." -> "I could not care less."
$ws.Range("B5").Value = "I could not care less."

# Update the active selection to B5
$ws.Range("B5").Select()
